$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for each value in a data row (G stays 1, H stays blank - untouched)
$cols = @('B','C','D','E','F','I','J','K','L','M','N')

# New vm_pu values for rows 2..25 (bus index 0..23), columns B,C,D,E,F,I,J,K,L,M,N
$data = @(
    @(1.02, 1.129247026424492, 1.130722603127266, 1.127063990637531, 1.140981927338763, 1.02359499962809, 1.133896276542335, 1.133270882220993, 1.129620971002904, 1.143506125943593, 1.135506539314746),
    @(1.02, 1.133225707344483, 1.13445635675796, 1.130691390990534, 1.144817092886285, 1.023504579208684, 1.13754103835214, 1.136825465177157, 1.133068768889232, 1.147163744512368, 1.139156477104381),
    @(1.02, 1.135778047712897, 1.136850916933997, 1.13301740081249, 1.147277086463472, 1.023443335729026, 1.139877754619948, 1.139104018246758, 1.135278475611616, 1.149508788926257, 1.141496511777315),
    @(1.02, 1.136845932583763, 1.137852628575423, 1.133990353537638, 1.148306265106009, 1.023416932628352, 1.140855084841424, 1.140056938218033, 1.136202508391054, 1.150489622383769, 1.142475229919806),
    @(1.02, 1.137024940099218, 1.138020534239541, 1.134153433827011, 1.148478780285696, 1.023412460900047, 1.141018892970115, 1.140216650092962, 1.136357373131804, 1.150654018841646, 1.142639270674829),
    @(1.02, 1.135792336721854, 1.136864321128302, 1.133030420500042, 1.147290857830439, 1.02344298551034, 1.139890833281398, 1.13911677057864, 1.135290841737638, 1.149521914390006, 1.141509609011964),
    @(1.02, 1.130596332553638, 1.131988982407889, 1.128294369322829, 1.142282620110822, 1.023565007130014, 1.135132633570558, 1.134476721655371, 1.130790667433168, 1.144746826078858, 1.136744652111761),
    @(1.02, 1.121262506868626, 1.123226112099651, 1.119779189702401, 1.133283969947179, 1.023759139049229, 1.126574254515701, 1.126128148900325, 1.122690686595267, 1.136158758830101, 1.128174119176986),
    @(1.02, 1.114908459795777, 1.11725740358316, 1.113977409134183, 1.127156767337608, 1.023874603872166, 1.120740681761561, 1.120435767421265, 1.117165735317105, 1.130305446719623, 1.122332262080415),
    @(1.02, 1.112123128571422, 1.114640213447459, 1.111432994774481, 1.124470582955494, 1.023921309042974, 1.118181747421177, 1.117938338832568, 1.114741274509411, 1.127737978669278, 1.119769693759755),
    @(1.02, 1.111083180531546, 1.113662928943056, 1.110482822726677, 1.123467611216251, 1.023938164268118, 1.117226064858701, 1.117005562577924, 1.113835678015827, 1.126779128449706, 1.118812654018421),
    @(1.02, 1.111306499163177, 1.113872796351518, 1.110686870689612, 1.123682991404904, 1.023934571058996, 1.117431300211458, 1.117205881638883, 1.114030163379013, 1.126985043176906, 1.11901818082892),
    @(1.02, 1.112037276777107, 1.114559537098778, 1.11135455787676, 1.124387784538379, 1.023922712353274, 1.118102857435525, 1.117861341029695, 1.114666521803922, 1.127658826799748, 1.119690691741274),
    @(1.02, 1.112486816273384, 1.114981972111548, 1.111765264183705, 1.124821334163844, 1.023915340507297, 1.118515932557203, 1.118264505889456, 1.115057928211831, 1.128073273966436, 1.120104353477005),
    @(1.02, 1.11509257622649, 1.117430389432967, 1.114145576023771, 1.127334324005294, 1.023871434968932, 1.120909796084229, 1.120600807997515, 1.117325943637882, 1.130475127691015, 1.122501616564819),
    @(1.02, 1.116717830811897, 1.118957302163238, 1.115629904598439, 1.128901640234791, 1.023843013862359, 1.122402420768214, 1.122057428110147, 1.118739858982833, 1.131972768117045, 1.123996360947062),
    @(1.02, 1.1176625521672, 1.1198447854996, 1.116492596651915, 1.129812655341854, 1.023826118601775, 1.123269879347639, 1.122903921094693, 1.11956148788879, 1.132843153594855, 1.124865051417169),
    @(1.02, 1.11798413172366, 1.120146868805908, 1.116786234623622, 1.130122756596417, 1.023820303822167, 1.123565130407845, 1.123192029264545, 1.11984112549347, 1.133139402869379, 1.125160721767746),
    @(1.02, 1.116543795541704, 1.118793805201203, 1.115470971575038, 1.128733812004513, 1.023846096030143, 1.122242605172769, 1.121901471782342, 1.118588479788393, 1.13181241424466, 1.123836318395141),
    @(1.02, 1.111822230749608, 1.114357452952329, 1.111158082263437, 1.124180386101842, 1.02392621804881, 1.117905245458828, 1.117668467757861, 1.114479271137417, 1.127460559144061, 1.119492799132904),
    @(1.02, 1.10882250157018, 1.111538262170433, 1.108416985200755, 1.121287238534612, 1.023973741382756, 1.115148083048271, 1.11497726954908, 1.111866352336078, 1.124694295365068, 1.116731721235472),
    @(1.02, 1.110415746786783, 1.113035680293403, 1.109872957988028, 1.122823897189455, 1.023948818288664, 1.116612638168607, 1.116406820780211, 1.113254362070879, 1.1261636740746, 1.118198356192066),
    @(1.02, 1.116622444653475, 1.11886769212567, 1.115542796106206, 1.128809656150222, 1.023844704313072, 1.12231482872387, 1.12197195138787, 1.118656891046031, 1.131884881015936, 1.123908644511969),
    @(1.02, 1.123697798587808, 1.125513024928247, 1.122001772998814, 1.135632050216929, 1.02371141742819, 1.128808519098924, 1.12830795879785, 1.124805950663286, 1.138400689581398, 1.130411556672191)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
